$d = $word.ActiveDocument

# Find the "Ver no Jupiter ..." paragraph - this is the first of three
# trailing paragraphs (itself, the copyright line, and a blank paragraph
# right after it) that need to be removed, leaving the blank "Normal"
# paragraph that precedes it directly followed by the page-break paragraph.

$verIndex = 0
For ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    If ($d.Paragraphs($i).Range.Text -eq "Ver no Jupiter Salvar em pdf Salvar em docx`r") {
        $verIndex = $i
        Break
    }
}

# Delete the blank paragraph after the copyright line, then the copyright
# line itself, then the "Ver no Jupiter ..." line - highest index first so
# the remaining indices don't shift underneath us.
$d.Paragraphs($verIndex + 2).Range.Delete()
$d.Paragraphs($verIndex + 1).Range.Delete()
$d.Paragraphs($verIndex).Range.Delete()
